# Weekly Fruit/Vegetable price update: refresh the "Guayaba" price series.
# Two brand-new price observations are inserted at the top (rows 8-9),
# which pushes every subsequent observation down by two rows; the two
# oldest observations that "fall off the bottom" are re-appended as new
# rows 57-58 (the table is a fixed-size rolling window of observations).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New values, per target row, for the columns that actually change ----
# (Fecha / Calidad / Volumen / Precio minimo / Precio maximo / Precio
#  promedio ponderado / Precio $/Kg)

$DVal = @{
    8 = 45051;  9 = 45051;  10 = 44798; 11 = 44798; 12 = 44386; 13 = 44386
    14 = 44344; 15 = 44344; 16 = 44260; 17 = 44764; 18 = 44778; 19 = 44778
    20 = 44750; 21 = 44379; 22 = 44379; 23 = 44407; 24 = 44403; 25 = 44403
    26 = 44767; 27 = 44767; 28 = 45030; 29 = 44372; 30 = 44372; 31 = 44725
    32 = 44725; 33 = 44722; 34 = 44722; 35 = 44694; 36 = 44694; 37 = 44309
    38 = 44890; 39 = 44473; 40 = 44316; 41 = 44690; 42 = 44350; 43 = 44414
    44 = 44351; 45 = 44351; 46 = 44389; 47 = 44389; 48 = 44687; 49 = 44326
    50 = 44358; 51 = 44358; 52 = 44715; 53 = 44715; 54 = 44498; 55 = 44799
    56 = 44330; 57 = 44330; 58 = 44883
}

$LVal = @{
    8 = "Primera";  9 = "Segunda";  10 = "Primera"; 11 = "Segunda"
    12 = "Primera"; 13 = "Segunda"; 14 = "Primera"; 15 = "Segunda"
    16 = "Primera"; 17 = "Primera"; 18 = "Primera"; 19 = "Segunda"
    20 = "Primera"; 21 = "Primera"; 22 = "Segunda"; 23 = "Primera"
    24 = "Primera"; 25 = "Segunda"; 26 = "Primera"; 27 = "Segunda"
    28 = "Primera"; 29 = "Primera"; 30 = "Segunda"; 31 = "Primera"
    32 = "Segunda"; 33 = "Primera"; 34 = "Segunda"; 35 = "Primera"
    36 = "Segunda"; 37 = "Primera"; 38 = "Segunda"; 39 = "Primera"
    40 = "Primera"; 41 = "Primera"; 42 = "Primera"; 43 = "Primera"
    44 = "Primera"; 45 = "Segunda"; 46 = "Primera"; 47 = "Segunda"
    48 = "Primera"; 49 = "Primera"; 50 = "Primera"; 51 = "Segunda"
    52 = "Primera"; 53 = "Segunda"; 54 = "Segunda"; 55 = "Primera"
    56 = "Primera"; 57 = "Segunda"; 58 = "Primera"
}

$MVal = @{
    8 = 120; 9 = 50;  10 = 100; 11 = 130; 12 = 160; 13 = 200; 14 = 140
    15 = 120; 16 = 100; 17 = 200; 18 = 200; 19 = 140; 20 = 200; 21 = 150
    22 = 140; 23 = 200; 24 = 100; 25 = 120; 26 = 200; 27 = 200; 28 = 100
    29 = 900; 30 = 900; 31 = 140; 32 = 160; 33 = 140; 34 = 200; 35 = 120
    36 = 140; 37 = 160; 38 = 160; 39 = 160; 40 = 140; 41 = 100; 42 = 140
    43 = 160; 44 = 100; 45 = 100; 46 = 140; 47 = 120; 48 = 120; 49 = 160
    50 = 200; 51 = 200; 52 = 120; 53 = 160; 54 = 100; 55 = 200; 56 = 200
    57 = 100; 58 = 100
}

$NVal = @{
    8 = 1600; 9 = 1400; 10 = 700;  11 = 500;  12 = 700;  13 = 600
    14 = 1000; 15 = 800;  16 = 1900; 17 = 500;  18 = 700;  19 = 500
    20 = 700;  21 = 700;  22 = 500;  23 = 600;  24 = 1200; 25 = 950
    26 = 800;  27 = 600;  28 = 1900; 29 = 750;  30 = 600;  31 = 700
    32 = 500;  33 = 800;  34 = 700;  35 = 1400; 36 = 1100; 37 = 1400
    38 = 1000; 39 = 1500; 40 = 1100; 41 = 1600; 42 = 750;  43 = 1300
    44 = 700;  45 = 600;  46 = 750;  47 = 600;  48 = 1300; 49 = 600
    50 = 700;  51 = 600;  52 = 800;  53 = 600;  54 = 1200; 55 = 600
    56 = 1200; 57 = 1000; 58 = 700
}

$OVal = @{
    8 = 1700; 9 = 1500; 10 = 800;  11 = 600;  12 = 750;  13 = 650
    14 = 1200; 15 = 850;  16 = 2000; 17 = 600;  18 = 800;  19 = 600
    20 = 800;  21 = 800;  22 = 600;  23 = 650;  24 = 1300; 25 = 1000
    26 = 900;  27 = 700;  28 = 2000; 29 = 800;  30 = 650;  31 = 800
    32 = 600;  33 = 900;  34 = 800;  35 = 1500; 36 = 1200; 37 = 1500
    38 = 1200; 39 = 1600; 40 = 1200; 41 = 1700; 42 = 800;  43 = 1400
    44 = 800;  45 = 700;  46 = 800;  47 = 700;  48 = 1400; 49 = 700
    50 = 800;  51 = 650;  52 = 900;  53 = 700;  54 = 1300; 55 = 700
    56 = 1300; 57 = 1100; 58 = 750
}

$PVal = @{
    8 = 1675; 9 = 1460; 10 = 750;  11 = 550;  12 = 725;  13 = 625
    14 = 1100; 15 = 825;  16 = 1950; 17 = 550;  18 = 750;  19 = 550
    20 = 750;  21 = 747;  22 = 543;  23 = 625;  24 = 1250; 25 = 975
    26 = 850;  27 = 650;  28 = 1950; 29 = 772;  30 = 628;  31 = 750
    32 = 550;  33 = 850;  34 = 750;  35 = 1450; 36 = 1150; 37 = 1450
    38 = 1100; 39 = 1550; 40 = 1150; 41 = 1650; 42 = 775;  43 = 1350
    44 = 750;  45 = 650;  46 = 775;  47 = 650;  48 = 1350; 49 = 650
    50 = 750;  51 = 625;  52 = 850;  53 = 650;  54 = 1250; 55 = 650
    56 = 1250; 57 = 1050; 58 = 725
}

$SVal = @{
    8 = 1675; 9 = 1460; 10 = 750;  11 = 550;  12 = 725;  13 = 625
    14 = 1100; 15 = 825;  16 = 1950; 17 = 550;  18 = 750;  19 = 550
    20 = 750;  21 = 747;  22 = 543;  23 = 625;  24 = 1250; 25 = 975
    26 = 850;  27 = 650;  28 = 1950; 29 = 772;  30 = 628;  31 = 750
    32 = 550;  33 = 850;  34 = 750;  35 = 1450; 36 = 1150; 37 = 1450
    38 = 1100; 39 = 1550; 40 = 1150; 41 = 1650; 42 = 775;  43 = 1350
    44 = 750;  45 = 650;  46 = 775;  47 = 650;  48 = 1350; 49 = 650
    50 = 750;  51 = 625;  52 = 850;  53 = 650;  54 = 1250; 55 = 650
    56 = 1250; 57 = 1050; 58 = 725
}

# ---- Columns that are identical on every data row of this sheet ----
$ConstA = 1
$ConstB = "Agrícola del Norte S.A. de Arica"
$ConstC = "Arica y Parinacota"
$ConstE = 15
$ConstF = "Fruta"
$ConstG = 100108
$ConstH = "Tropicales y subtropicales"
$ConstI = 100108001
$ConstJ = "Guayaba"
$ConstK = "Sin especificar"
$ConstQ = "$/kilo (en caja de 10 kilos )"
$ConstR = "Región de Arica y Parinacota"
$ConstT = 1

# ---- Apply the updated values to the existing rows (8-56) ----
for ($r = 8; $r -le 56; $r++) {
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 4).Value = $DVal[$r]
    $ws.Cells.Item($r, 12).Value = $LVal[$r]
    $ws.Cells.Item($r, 13).Value = $MVal[$r]
    $ws.Cells.Item($r, 14).Value = $NVal[$r]
    $ws.Cells.Item($r, 15).Value = $OVal[$r]
    $ws.Cells.Item($r, 16).Value = $PVal[$r]
    $ws.Cells.Item($r, 19).Value = $SVal[$r]
}

# ---- Append the two new rows (57-58) that extend the table ----
for ($r = 57; $r -le 58; $r++) {
    $ws.Cells.Item($r, 1).Value  = $ConstA
    $ws.Cells.Item($r, 2).Value  = $ConstB
    $ws.Cells.Item($r, 3).Value  = $ConstC
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 4).Value  = $DVal[$r]
    $ws.Cells.Item($r, 5).Value  = $ConstE
    $ws.Cells.Item($r, 6).Value  = $ConstF
    $ws.Cells.Item($r, 7).Value  = $ConstG
    $ws.Cells.Item($r, 8).Value  = $ConstH
    $ws.Cells.Item($r, 9).Value  = $ConstI
    $ws.Cells.Item($r, 10).Value = $ConstJ
    $ws.Cells.Item($r, 11).Value = $ConstK
    $ws.Cells.Item($r, 12).Value = $LVal[$r]
    $ws.Cells.Item($r, 13).Value = $MVal[$r]
    $ws.Cells.Item($r, 14).Value = $NVal[$r]
    $ws.Cells.Item($r, 15).Value = $OVal[$r]
    $ws.Cells.Item($r, 16).Value = $PVal[$r]
    $ws.Cells.Item($r, 17).Value = $ConstQ
    $ws.Cells.Item($r, 18).Value = $ConstR
    $ws.Cells.Item($r, 19).Value = $SVal[$r]
    $ws.Cells.Item($r, 20).Value = $ConstT
}
